$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# ---------------------------------------------------------------------------
# New activity-log entries for rows 24-26 (last 4 digits / date / start /
# end time / description) that were previously blank placeholder rows.
# ---------------------------------------------------------------------------

$ws.Range("B24").Value = 6977
$ws.Range("C24").Value = 43924
$ws.Range("D24").Value = 0.88611111111111107
$ws.Range("E24").Value = 0.89930555555555547
$ws.Range("G24").Value = "Fixed folder structure and files as they were incorrect. Fixed project file settings"

$ws.Range("B25").Value = 6977
$ws.Range("C25").Value = 43924
$ws.Range("D25").Value = 0.89930555555555547
$ws.Range("E25").Value = 0.91319444444444453
$ws.Range("G25").Value = "Ran scripts for Test Benches and ConfigExU. Fixed changes along the way (but still does not run successfully)"

$ws.Range("B26").Value = 6977
$ws.Range("C26").Value = 43924
$ws.Range("D26").Value = 0.91319444444444453
$ws.Range("E26").Value = 0.92708333333333337
$ws.Range("G26").Value = "Fixed ArithUnit.vhd and Adder.vhd so they compile on ModelSim and Quartus."

# ---------------------------------------------------------------------------
# Move the cursor/selection to match where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Range("C26").Select()

# ---------------------------------------------------------------------------
# Resize/reposition the workbook window (maximised on save).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 15840
